$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.009.93"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.298.35"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "2.289.81"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "42.923.14"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  +4.16%  "
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +7.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0839"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.127"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  +11.75%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.04%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").Value = "BinanceUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
